$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value. All target cells are plain text
# (coin names, URLs, "." -separated price strings, padded percent strings),
# so force text format before assigning to stop Excel from auto-coercing
# values like "36.584.05" or "0.602" into numbers/dates.
$updates = [ordered]@{
    "D2" = "36.584.05"
    "E2" = "  -1.09%  "
    "D3" = "2.026.01"
    "E3" = "  +1.49%  "
    "E4" = "  +0.07%  "
    "E5" = "  -8.87%  "
    "D6" = "0.602"
    "E6" = "  -2.21%  "
    "E7" = "  +0.02%  "
    "D8" = "55.44"
    "E8" = "  +0.42%  "
    "E9" = "  -0.76%  "
    "D10" = "57.23"
    "E10" = "  +2.78%  "
    "D11" = "0.0752"
    "E11" = "  -1.64%  "
    "E12" = "  -1.19%  "
    "D13" = "2.324.68"
    "E13" = "  +1.51%  "
    "D14" = "14.38"
    "E14" = "  +1.69%  "
    "D15" = "20.26"
    "E15" = "  -4.49%  "
    "E16" = "  -3.25%  "
    "D17" = "5.11"
    "E17" = "  -0.92%  "
    "D18" = "2.038.08"
    "E18" = "  +2.80%  "
    "D19" = "36.735.76"
    "E19" = "  -0.42%  "
    "D20" = "67.83"
    "E20" = "  -3.76%  "
    "D21" = "5.57"
    "E21" = "  +10.34%  "
    "E22" = "  -2.78%  "
    "D23" = "221.65"
    "E23" = "  -5.20%  "
    "D24" = "0.999"
    "E24" = "  -0.21%  "
    "D25" = "2.39"
    "E25" = "  +1.05%  "
    "D26" = "2.38"
    "E26" = "  -5.38%  "
    "D27" = "162.96"
    "E27" = "  -0.67%  "
    "D28" = "0.132"
    "E28" = "  +6.59%  "
    "E29" = "  -2.03%  "
    "D30" = "18.98"
    "E30" = "  -2.15%  "
    "E31" = "  +2.05%  "
    "E32" = "  -1.05%  "
    "E33" = "  -3.31%  "
    "E34" = "  -4.43%  "
    "D35" = "2.46"
    "E35" = "  +5.29%  "
    "D36" = "4.27"
    "E36" = "  -2.32%  "
    "E37" = "  +0.25%  "
    "B38" = "WEMIXToken"
    "C38" = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
    "D38" = "1.77"
    "E38" = "  -2.73%  "
    "B39" = "RenderToken"
    "C39" = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
    "D39" = "3.29"
    "E39" = "  -6.38%  "
    "E40" = "  +5.51%  "
    "D41" = "2.90"
    "E41" = "  -3.67%  "
    "D42" = "0.0945"
    "E42" = "  +3.58%  "
    "D43" = "1.462.67"
    "E43" = "  +1.23%  "
    "D44" = "93.71"
    "E44" = "  +5.77%  "
    "E45" = "  -2.12%  "
    "B46" = "TrustWalletToken"
    "C46" = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
    "D46" = "1.11"
    "E46" = "  -4.63%  "
    "B47" = "InjectiveProtocol"
    "C47" = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
    "D47" = "15.73"
    "E47" = "  +1.61%  "
    "B48" = "FTXToken"
    "C48" = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
    "D48" = "4.06"
    "E48" = "  +38.17%  "
    "E49" = "  -0.92%  "
    "D50" = "2.89"
    "E50" = "  -0.04%  "
    "E51" = "  +0.64%  "
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
}
